$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. Insert a row at
# position 72 (shifting the existing rows 72-175 down to 73-176) and
# populate it with the new record's values.
$ws.Rows(72).Insert()

$ws.Cells.Item(72, 1).Value = 7
$ws.Cells.Item(72, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(72, 3).Value = "Ñuble"
$ws.Cells.Item(72, 4).Value = 44571
$ws.Cells.Item(72, 5).Value = 16
$ws.Cells.Item(72, 6).Value = 100112032
$ws.Cells.Item(72, 7).Value = "Zapallo italiano"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 120
$ws.Cells.Item(72, 11).Value = 6500
$ws.Cells.Item(72, 12).Value = 7000
$ws.Cells.Item(72, 13).Value = 6750
$ws.Cells.Item(72, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(72, 15).Value = "Región del Maule"
$ws.Cells.Item(72, 16).Value = 112
$ws.Cells.Item(72, 17).Value = 60
$ws.Cells.Item(72, 18).Value = "Hortaliza"
